# Fix header label "birth_date" -> "birth date" across the tutorial sheets.
$wb = $excel.ActiveWorkbook

$sheetsWithColumnC = @("Animal", "NamedThing", "Animal1", "NamedThing1")
foreach ($name in $sheetsWithColumnC) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = "birth date"
}

$sheetsWithColumnB = @("Person", "Person1")
foreach ($name in $sheetsWithColumnB) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "birth date"
}
